# Apply updated vulnerability-detection metrics (SWC-112 -> DASP-2 linkage)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: access_control
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 6
$ws.Range("E2").Value = 0.6470588235294118
$ws.Range("F2").Value = 0.6470588235294118

# Row 4: denial_service
$ws.Range("D4").Value = 2

# Row 8: front_running
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 0.05555555555555555
$ws.Range("F8").Value = 0.5

# Row 11: Other
$ws.Range("D11").Value = 18
